$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 98.912777
$ws.Range("H2").Value = 296.738331
$ws.Range("I2").Value = 0.8120825131376513
$ws.Range("J2").Value = 0.8120825131376513
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 114.155417
$ws.Range("N2").Value = 342.466251
$ws.Range("O2").Value = 0.6835107367845005
$ws.Range("P2").Value = 0.6835107367845005
$ws.Range("Q2").Value = 11291.42930506301
$ws.Range("R2").Value = 101622.8637455671
$ws.Range("S2").Value = 0.5550671168845248
$ws.Range("T2").Value = 0.5550671168845248

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 98.912777
$ws.Range("H3").Value = 296.738331
$ws.Range("I3").Value = 0.8120825131376513
$ws.Range("J3").Value = 0.8120825131376513
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.924535
$ws.Range("N3").Value = 107.773605
$ws.Range("O3").Value = 0.2150997826628812
$ws.Range("P3").Value = 0.2150997826628812
$ws.Range("Q3").Value = 3553.395519283695
$ws.Range("R3").Value = 31980.55967355326
$ws.Range("S3").Value = 0.1746787720802352
$ws.Range("T3").Value = 0.1746787720802352

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 98.912777
$ws.Range("H4").Value = 296.738331
$ws.Range("I4").Value = 0.8120825131376513
$ws.Range("J4").Value = 0.8120825131376513
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.93339666666667
$ws.Range("N4").Value = 50.80019
$ws.Range("O4").Value = 0.1013894805526183
$ws.Range("P4").Value = 0.1013894805526183
$ws.Range("Q4").Value = 1674.929288342543
$ws.Range("R4").Value = 15074.36359508289
$ws.Range("S4").Value = 0.08233662417289134
$ws.Range("T4").Value = 0.08233662417289134

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 17.04862266666667
$ws.Range("H5").Value = 51.14586800000001
$ws.Range("I5").Value = 0.1399706767982279
$ws.Range("J5").Value = 0.1399706767982279
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 114.155417
$ws.Range("N5").Value = 342.466251
$ws.Range("O5").Value = 0.6835107367845005
$ws.Range("P5").Value = 0.6835107367845005
$ws.Range("Q5").Value = 1946.192629788986
$ws.Range("R5").Value = 17515.73366810087
$ws.Range("S5").Value = 0.09567146042658196
$ws.Range("T5").Value = 0.09567146042658196

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.04862266666667
$ws.Range("H6").Value = 51.14586800000001
$ws.Range("I6").Value = 0.1399706767982279
$ws.Range("J6").Value = 0.1399706767982279
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.924535
$ws.Range("N6").Value = 107.773605
$ws.Range("O6").Value = 0.2150997826628812
$ws.Range("P6").Value = 0.2150997826628812
$ws.Range("Q6").Value = 612.4638416904601
$ws.Range("R6").Value = 5512.174575214141
$ws.Range("S6").Value = 0.03010766215847522
$ws.Range("T6").Value = 0.03010766215847522

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.04862266666667
$ws.Range("H7").Value = 51.14586800000001
$ws.Range("I7").Value = 0.1399706767982279
$ws.Range("J7").Value = 0.1399706767982279
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.93339666666667
$ws.Range("N7").Value = 50.80019
$ws.Range("O7").Value = 0.1013894805526183
$ws.Range("P7").Value = 0.1013894805526183
$ws.Range("Q7").Value = 288.6910902349912
$ws.Range("R7").Value = 2598.21981211492
$ws.Range("S7").Value = 0.01419155421317076
$ws.Range("T7").Value = 0.01419155421317076

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.839988000000001
$ws.Range("H8").Value = 17.519964
$ws.Range("I8").Value = 0.0479468100641207
$ws.Range("J8").Value = 0.04794681006412069
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 114.155417
$ws.Range("N8").Value = 342.466251
$ws.Range("O8").Value = 0.6835107367845005
$ws.Range("P8").Value = 0.6835107367845005
$ws.Range("Q8").Value = 666.6662654149961
$ws.Range("R8").Value = 5999.996388734965
$ws.Range("S8").Value = 0.03277215947339365
$ws.Range("T8").Value = 0.03277215947339364

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.839988000000001
$ws.Range("H9").Value = 17.519964
$ws.Range("I9").Value = 0.0479468100641207
$ws.Range("J9").Value = 0.04794681006412069
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 35.924535
$ws.Range("N9").Value = 107.773605
$ws.Range("O9").Value = 0.2150997826628812
$ws.Range("P9").Value = 0.2150997826628812
$ws.Range("Q9").Value = 209.79885330558
$ws.Range("R9").Value = 1888.18967975022
$ws.Range("S9").Value = 0.01031334842417081
$ws.Range("T9").Value = 0.01031334842417081

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.839988000000001
$ws.Range("H10").Value = 17.519964
$ws.Range("I10").Value = 0.0479468100641207
$ws.Range("J10").Value = 0.04794681006412069
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.93339666666667
$ws.Range("N10").Value = 50.80019
$ws.Range("O10").Value = 0.1013894805526183
$ws.Range("P10").Value = 0.1013894805526183
$ws.Range("Q10").Value = 98.89083333257335
$ws.Range("R10").Value = 890.0174999931601
$ws.Range("S10").Value = 0.004861302166556251
$ws.Range("T10").Value = 0.00486130216655625
